$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that became empty inline strings
$ws.Range("I2").Value = ""

$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("J3").Value = ""

$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("J4").Value = ""

$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("J5").Value = ""

$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("J6").Value = ""

$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("J8").Value = ""

$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("J9").Value = ""

$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("J10").Value = ""

$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("J11").Value = ""

$ws.Range("D12").Value = ""

$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("J13").Value = ""
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("J14").Value = ""

$ws.Range("C24").Value = ""
$ws.Range("C25").Value = ""
